# Append weeks 19-22 of the power rankings to the bottom of the sheet.
# Each week repeats the same 12-team ordering (ranks 1-12) that week 18
# already uses, starting at row 218 and ending at row 265.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$teams = @(
    "Baby Back Gibbs",
    "Magic Mikaela",
    "Sith Happens",
    "Apex Predators",
    "Mighty Rubber Ducks",
    "Compile and Conquer",
    "Drafted by AI",
    "Aida's Astounding Team",
    "Bring the heat",
    "Kelly's Deluxe Team",
    "Boomer Sooners",
    "Kuppenheimer"
)

$row = 218
for ($week = 19; $week -le 22; $week++) {
    for ($i = 0; $i -lt $teams.Length; $i++) {
        $ws.Cells.Item($row, 1).Value = $teams[$i]
        $ws.Cells.Item($row, 2).Value = $week
        $ws.Cells.Item($row, 3).Value = ($i + 1)
        $row++
    }
}
